$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect before editing, then restore protection after.
$ws.Unprotect()

# Update the confidential disclosure date string (A59)
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-10 for illustrative purposes only and are subject to change."
$ws.Range("A59").Value2 = $newText
$ws.Rows.Item(59).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-56
$ws.Range("D2").Value2 = [double]"0.01467611323531555"
$ws.Range("E2").Value2 = [double]"0.003010679391047377"
$ws.Range("D3").Value2 = [double]"0.05065613133174238"
$ws.Range("E3").Value2 = [double]"0.02087682672233826"
$ws.Range("D4").Value2 = [double]"0.01383940668619249"
$ws.Range("E4").Value2 = [double]"0.03725932040784485"
$ws.Range("D5").Value2 = [double]"0.00949531323294419"
$ws.Range("E5").Value2 = [double]"0.002081165452653577"
$ws.Range("D6").Value2 = [double]"0.01533787036723692"
$ws.Range("E6").Value2 = [double]"0.01026694045174548"
$ws.Range("D7").Value2 = [double]"0.01956101837975472"
$ws.Range("E7").Value2 = [double]"0.007418947993174552"
$ws.Range("D8").Value2 = [double]"0.005120126581121568"
$ws.Range("E8").Value2 = [double]"-0.01194647304656171"
$ws.Range("D9").Value2 = [double]"0.006964011922991131"
$ws.Range("E9").Value2 = [double]"0.001596169193934571"
$ws.Range("D10").Value2 = [double]"0.01425176981127694"
$ws.Range("E10").Value2 = [double]"-0.006844106463878341"
$ws.Range("D11").Value2 = [double]"0.007984807498827289"
$ws.Range("E11").Value2 = [double]"0.002706883217324041"
$ws.Range("D12").Value2 = [double]"0.01524709181330594"
$ws.Range("E12").Value2 = [double]"-0.01822600243013361"
$ws.Range("D13").Value2 = [double]"0.003290722579998137"
$ws.Range("E13").Value2 = [double]"-0.02036124794745486"
$ws.Range("D14").Value2 = [double]"0.006049124660074455"
$ws.Range("E14").Value2 = [double]"0.01122964626614276"
$ws.Range("D15").Value2 = [double]"0.01459021325673188"
$ws.Range("E15").Value2 = [double]"-0.0155885602062108"
$ws.Range("D16").Value2 = [double]"0.01085730031250054"
$ws.Range("E16").Value2 = [double]"-0.0159969285897108"
$ws.Range("D17").Value2 = [double]"0.02230234544611545"
$ws.Range("E17").Value2 = [double]"0.01675908043112573"
$ws.Range("D18").Value2 = [double]"0.008481310610123279"
$ws.Range("E18").Value2 = [double]"0.03042813455657489"
$ws.Range("D19").Value2 = [double]"0.01636138312427917"
$ws.Range("E19").Value2 = [double]"0.008998127906274522"
$ws.Range("D20").Value2 = [double]"0.01395643074040965"
$ws.Range("E20").Value2 = [double]"0.03336283185840716"
$ws.Range("D21").Value2 = [double]"0.006628501798156417"
$ws.Range("E21").Value2 = [double]"0.005860048259220818"
$ws.Range("D22").Value2 = [double]"0.01468883458368955"
$ws.Range("E22").Value2 = [double]"0.006474396703943563"
$ws.Range("D23").Value2 = [double]"0.01936806763706264"
$ws.Range("E23").Value2 = [double]"0.00154958677685979"
$ws.Range("D24").Value2 = [double]"0.01003334598881065"
$ws.Range("E24").Value2 = [double]"-0.01481481481481473"
$ws.Range("D25").Value2 = [double]"0.02129998347459794"
$ws.Range("E25").Value2 = [double]"-0.0001507613447913059"
$ws.Range("D26").Value2 = [double]"0.0130370971809739"
$ws.Range("E26").Value2 = [double]"0.01064358897841422"
$ws.Range("D27").Value2 = [double]"0.02153671788717918"
$ws.Range("E27").Value2 = [double]"0.01053478212871961"
$ws.Range("D28").Value2 = [double]"0.05495560743450419"
$ws.Range("E28").Value2 = [double]"-0.008023283253362656"
$ws.Range("D29").Value2 = [double]"0.01993657605032021"
$ws.Range("E29").Value2 = [double]"0.00654817586529477"
$ws.Range("D30").Value2 = [double]"0.03055223249783304"
$ws.Range("E30").Value2 = [double]"0.01305736346363751"
$ws.Range("D31").Value2 = [double]"0.0151820338496554"
$ws.Range("E31").Value2 = [double]"0.01000827752276345"
$ws.Range("D32").Value2 = [double]"0.01306970335544707"
$ws.Range("E32").Value2 = [double]"0.007361522578328339"
$ws.Range("D33").Value2 = [double]"0.01748728420949664"
$ws.Range("E33").Value2 = [double]"0.007882023900330548"
$ws.Range("D34").Value2 = [double]"0.04461006350052478"
$ws.Range("E34").Value2 = [double]"0.01129180959658482"
$ws.Range("D35").Value2 = [double]"0.0107415113406498"
$ws.Range("E35").Value2 = [double]"0.003104518799586176"
$ws.Range("D36").Value2 = [double]"0.009855957283935308"
$ws.Range("E36").Value2 = [double]"0.01964912280701747"
$ws.Range("D37").Value2 = [double]"0.0107730059409932"
$ws.Range("E37").Value2 = [double]"-0.005159071367153989"
$ws.Range("D38").Value2 = [double]"0.007169653136998052"
$ws.Range("E38").Value2 = [double]"-0.005943152454780387"
$ws.Range("D39").Value2 = [double]"0.01215246943603422"
$ws.Range("E39").Value2 = [double]"-0.0148865784499056"
$ws.Range("D40").Value2 = [double]"0.01768057459983948"
$ws.Range("E40").Value2 = [double]"0.01869677897075173"
$ws.Range("D41").Value2 = [double]"0.01717869888024962"
$ws.Range("E41").Value2 = [double]"0.002800355166996837"
$ws.Range("D42").Value2 = [double]"0.03374831689150518"
$ws.Range("E42").Value2 = [double]"0.01756649191666915"
$ws.Range("D43").Value2 = [double]"0.0114763235945442"
$ws.Range("E43").Value2 = [double]"0.001656272364788824"
$ws.Range("D44").Value2 = [double]"0.02223645380186759"
$ws.Range("E44").Value2 = [double]"0.007059532521200174"
$ws.Range("D45").Value2 = [double]"0.01289231465057173"
$ws.Range("E45").Value2 = [double]"0.01021705860797972"
$ws.Range("D46").Value2 = [double]"0.00814533732940116"
$ws.Range("E46").Value2 = [double]"-0.0009211558800451547"
$ws.Range("D47").Value2 = [double]"0.01304154347749297"
$ws.Range("E47").Value2 = [double]"-0.03801406349882808"
$ws.Range("D48").Value2 = [double]"0.01056606789050155"
$ws.Range("E48").Value2 = [double]"-0.01396559886381565"
$ws.Range("D49").Value2 = [double]"0.01613523954298988"
$ws.Range("E49").Value2 = [double]"-0.004018646519852065"
$ws.Range("D50").Value2 = [double]"0.00850354209271862"
$ws.Range("E50").Value2 = [double]"0.009197530864197612"
$ws.Range("D51").Value2 = [double]"0.01063974055365778"
$ws.Range("E51").Value2 = [double]"-0.03131892368769296"
$ws.Range("D52").Value2 = [double]"0.008265325581297689"
$ws.Range("E52").Value2 = [double]"0.008465173618245281"
$ws.Range("D53").Value2 = [double]"0.009137293731981692"
$ws.Range("E53").Value2 = [double]"-0.00423418095801309"
$ws.Range("D54").Value2 = [double]"0.134122905022413"
$ws.Range("E54").Value2 = [double]"-9.853187506159156E-05"
$ws.Range("D55").Value2 = [double]"0.04412918418116385"
$ws.Range("E55").Value2 = [double]"0.004995836802664355"
$ws.Range("D56").Value2 = [double]"0.9999999999999999"
$ws.Range("E56").Value2 = [double]"0.004231348959430559"

$ws.Protect()
"done"
